$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.438.17"
$ws.Range("E2").Value = "  -1.95%  "
$ws.Range("D3").Value = "2.451.46"
$ws.Range("E3").Value = "  -2.28%  "
$ws.Range("E4").Value = "  +0.04%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "564.09"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -2.17%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "163.67"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -2.32%  "
$ws.Range("E7").Value = "  +0.02%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.507"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -1.40%  "
$ws.Range("E9").Value = "  -6.58%  "
$ws.Range("E10").Value = "  -2.05%  "
$ws.Range("E11").Value = "  -4.55%  "
$ws.Range("E12").Value = "  -2.68%  "
$ws.Range("D13").Value = "2.902.76"
$ws.Range("E13").Value = "  -1.67%  "
$ws.Range("D14").Value = "68.383.15"
$ws.Range("E14").Value = "  -1.82%  "
$ws.Range("E15").Value = "  -4.18%  "
$ws.Range("E16").Value = "  -5.15%  "
$ws.Range("D17").Value = "2.464.48"
$ws.Range("E17").Value = "  -0.79%  "
$ws.Range("E18").Value = "  -2.48%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "344.74"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -1.38%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "7.16"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -4.65%  "
$ws.Range("E21").Value = "  -2.23%  "
$ws.Range("E22").Value = "  -3.40%  "
$ws.Range("E23").Value = "  -0.09%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "68.15"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -5.34%  "
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "1.05"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +4.89%  "
$ws.Range("B27").Value = "WrappedeETH"
$ws.Range("C27").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D27").Value = "2.580.14"
$ws.Range("E27").Value = "  -0.40%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "8.27"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -6.26%  "
$ws.Range("E29").Value = "  -5.91%  "
$ws.Range("E30").Value = "  -6.89%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "437.81"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -5.17%  "
$ws.Range("E32").Value = "  -3.16%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +0.09%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.68"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -3.05%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "3.02"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +101.93%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "156.81"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -0.77%  "
$ws.Range("E37").Value = "  -0.40%  "
$ws.Range("E39").Value = "  -5.80%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "17.91"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -3.18%  "
$ws.Range("E41").Value = "  -3.88%  "
$ws.Range("E42").Value = "  -4.31%  "
$ws.Range("E43").Value = "  -4.39%  "
$ws.Range("E44").Value = "  +5.67%  "
$ws.Range("E45").Value = "  -5.45%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "135.26"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -4.46%  "
$ws.Range("E47").Value = "  -3.19%  "
$ws.Range("E48").Value = "  -2.47%  "
$ws.Range("E50").Value = "  -2.57%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.0917"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -1.45%  "
